$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.453.63"
$ws.Range("E2").Value = "  +1.11%  "
$ws.Range("D3").Value = "2.985.78"
$ws.Range("E3").Value = "  +2.80%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'381.61"
$ws.Range("E5").Value = "  +3.26%  "
$ws.Range("D6").Value = "'103.59"
$ws.Range("E6").Value = "  +1.37%  "
$ws.Range("D7").Value = "'0.548"
$ws.Range("E7").Value = "  +1.42%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("E9").Value = "  +2.28%  "
$ws.Range("E10").Value = "  +1.61%  "
$ws.Range("E11").Value = "  +0.29%  "
$ws.Range("D12").Value = "'0.0847"
$ws.Range("E12").Value = "  +1.83%  "
$ws.Range("D13").Value = "3.455.24"
$ws.Range("E13").Value = "  +3.10%  "
$ws.Range("D14").Value = "'18.40"
$ws.Range("E14").Value = "  +0.84%  "
$ws.Range("D15").Value = "'7.57"
$ws.Range("E15").Value = "  +3.07%  "
$ws.Range("D16").Value = "2.993.49"
$ws.Range("E16").Value = "  +3.47%  "
$ws.Range("D17").Value = "'0.970"
$ws.Range("E17").Value = "  +5.54%  "
$ws.Range("D18").Value = "51.449.52"
$ws.Range("E18").Value = "  +1.19%  "
$ws.Range("D19").Value = "'3.31"
$ws.Range("E19").Value = "  +3.74%  "
$ws.Range("D20").Value = "'7.43"
$ws.Range("E20").Value = "  +3.92%  "
$ws.Range("D21").Value = "'12.93"
$ws.Range("E21").Value = "  +0.57%  "
$ws.Range("D22").Value = "0.0₃0964"
$ws.Range("E22").Value = "  +2.60%  "
$ws.Range("D23").Value = "'69.07"
$ws.Range("E23").Value = "  +1.87%  "
$ws.Range("D24").Value = "'262.84"
$ws.Range("E24").Value = "  +1.93%  "
$ws.Range("D25").Value = "'2.94"
$ws.Range("E25").Value = "  +10.34%  "
$ws.Range("E26").Value = "  +17.25%  "
$ws.Range("D27").Value = "'7.72"
$ws.Range("E27").Value = "  +23.47%  "
$ws.Range("E28").Value = "  +15.31%  "
$ws.Range("E29").Value = "  +2.55%  "
$ws.Range("E30").Value = "  +1.92%  "
$ws.Range("E31").Value = "  -0.01%  "
$ws.Range("D32").Value = "'9.90"
$ws.Range("E32").Value = "  +0.58%  "
$ws.Range("D33").Value = "'34.77"
$ws.Range("E33").Value = "  +2.21%  "
$ws.Range("D34").Value = "'51.02"
$ws.Range("E34").Value = "  -0.60%  "
$ws.Range("E35").Value = "  -1.99%  "
$ws.Range("D36").Value = "'0.0453"
$ws.Range("E36").Value = "  +8.52%  "
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("E38").Value = "  +2.50%  "
$ws.Range("D39").Value = "'17.09"
$ws.Range("E39").Value = "  +0.79%  "
$ws.Range("E40").Value = "  +0.36%  "
$ws.Range("E41").Value = "  +0.55%  "
$ws.Range("E42").Value = "  +4.16%  "
$ws.Range("D43").Value = "'122.32"
$ws.Range("E43").Value = "  +2.82%  "
$ws.Range("D44").Value = "'21.90"
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("E45").Value = "  +18.57%  "
$ws.Range("E46").Value = "  -2.10%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "'3.27"
$ws.Range("E48").Value = "  +4.91%  "
$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").Value = "2.031.10"
$ws.Range("E49").Value = "  +0.81%  "
$ws.Range("D50").Value = "'0.0333"
$ws.Range("E50").Value = "  +8.54%  "
$ws.Range("D51").Value = "'58.23"
$ws.Range("E51").Value = "  +3.72%  "
